$wb = $excel.ActiveWorkbook

# Update "Impact Matrix" sheet (rows 2 and 3)
$wsImpact = $wb.Worksheets.Item("Impact Matrix")
$wsImpact.Range("B2").Value = "NDPS 2026-2030 Launch"
$wsImpact.Range("B3").Value = "IPS / Ethiopay Launch"

# Update "Events Metadata" sheet (rows 12 and 13)
$wsEvents = $wb.Worksheets.Item("Events Metadata")
$wsEvents.Range("B12").Value = "NDPS 2026-2030 Launch"
$wsEvents.Range("B13").Value = "IPS / Ethiopay Launch"
